# Add a new weekly price record for "Brocoli" (Vega Modelo de Temuco) as a
# new row 616. Inserting an entire row shifts the existing rows 616:653
# down to 617:654 (their data/formatting moves along unchanged), so all
# that remains is to populate the freshly-inserted row with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = -4121

# Insert a new blank row above the current row 616, pushing rows 616:653
# down to 617:654.
$ws.Rows("616:616").Insert($xlShiftDown)

# Populate the new row 616 with the new record.
$row = $ws.Range("A616:R616")
$row.Cells.Item(1, 1).Value = 10
$row.Cells.Item(1, 2).Value = "Vega Modelo de Temuco"
$row.Cells.Item(1, 3).Value = "La Araucanía"
$row.Cells.Item(1, 4).Value = 45041
$row.Cells.Item(1, 5).Value = 9
$row.Cells.Item(1, 6).Value = 100112023
$row.Cells.Item(1, 7).Value = "Brócoli"
$row.Cells.Item(1, 8).Value = "Sin especificar"
$row.Cells.Item(1, 9).Value = "Primera"
$row.Cells.Item(1, 10).Value = 1500
$row.Cells.Item(1, 11).Value = 1400
$row.Cells.Item(1, 12).Value = 1400
$row.Cells.Item(1, 13).Value = 1400
$row.Cells.Item(1, 14).Value = "$/unidad"
$row.Cells.Item(1, 15).Value = "Región Metropolitana"
$row.Cells.Item(1, 16).Value = 1400
$row.Cells.Item(1, 17).Value = 1
$row.Cells.Item(1, 18).Value = "Hortaliza"
